# Insert a new weekly price record as row 125 (Jengibre - Vega Central
# Mapocho de Santiago, Primera, fecha 45180) pushing the previously
# existing rows 125-144 down to 126-145.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift rows 125..144 down to 126..145, leaving a blank row 125.
$ws.Rows.Item(125).Insert()

# Populate the newly inserted row 125 with the new weekly record.
$ws.Cells.Item(125, 1).Value2 = 9
$ws.Cells.Item(125, 2).Value2 = "Vega Central Mapocho de Santiago"
$ws.Cells.Item(125, 3).Value2 = "Metropolitana"
$ws.Cells.Item(125, 4).Value2 = 45180
$ws.Cells.Item(125, 5).Value2 = 13
$ws.Cells.Item(125, 6).Value2 = 100114007
$ws.Cells.Item(125, 7).Value2 = "Jengibre"
$ws.Cells.Item(125, 8).Value2 = "Sin especificar"
$ws.Cells.Item(125, 9).Value2 = "Primera"
$ws.Cells.Item(125, 10).Value2 = 430
$ws.Cells.Item(125, 11).Value2 = 17000
$ws.Cells.Item(125, 12).Value2 = 18000
$ws.Cells.Item(125, 13).Value2 = 17500
$ws.Cells.Item(125, 14).Value2 = "`$/caja 13 kilos"
$ws.Cells.Item(125, 15).Value2 = "Perú"
$ws.Cells.Item(125, 16).Value2 = 1346
$ws.Cells.Item(125, 17).Value2 = 13
$ws.Cells.Item(125, 18).Value2 = "Hortaliza"
